# Connect To Camera Through USB Port and create Reload features
#
# This script mirrors what happened interactively in Excel:
#  1. New data was added to "Sheet6" (existing data was left untouched,
#     new cells E4/E6/I3/I5 were populated, and the selection moved).
#  2. A brand-new worksheet "Sheet7" was added at the end of the workbook,
#     filled with data, and made the active/selected sheet (tabSelected).
#  3. The workbook window was minimized, which shifts the active-tab
#     index forward by one (5 -> 6) to keep pointing at the (now)
#     newly-added last sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Update Sheet6 with additional values -----------------------------
$ws6 = $wb.Worksheets.Item("Sheet6")

$ws6.Range("I3").Value = 33
$ws6.Range("E4").Value = 63
$ws6.Range("I5").Value = 34
$ws6.Range("E6").Value = 90

# --- 2. Add the new Sheet7 worksheet at the end of the workbook ----------
$ws7 = $wb.Worksheets.Add()
$ws7.Name = "Sheet7"
$ws7.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
# Re-fetch by name: the reference returned by Add()/Move() can go stale
# (re-seats to whichever sheet now sits at the old index) once the sheet
# collection has been reordered.
$ws7 = $wb.Worksheets.Item("Sheet7")

$ws7.Range("C2").Value = 157
$ws7.Range("E2").Value = 34
$ws7.Range("G2").Value = 123

$ws7.Range("I3").Value = 129

$ws7.Range("C4").Value = 62
$ws7.Range("E4").Value = 157
$ws7.Range("J4").Value = 56
$ws7.Range("L4").Value = 31

$ws7.Range("F6").Value = 34
$ws7.Range("H6").Value = 94
$ws7.Range("L6").Value = 31

$ws7.Range("B7").Value = 157
$ws7.Range("K7").Value = 858

$ws7.Range("D8").Value = 94

$ws7.Range("G9").Value = 157

$ws7.Range("F10").Value = 157

# --- 3. Selections -- Sheet6 selection moves to K5, Sheet7's to K7 -------
$ws6 = $wb.Worksheets.Item("Sheet6")
$ws6.Range("K5").Select()
$ws7.Range("K7").Select()

# Sheet7 becomes the active / tab-selected sheet.
$ws7.Activate()

# --- 4. Window state: minimize the workbook window -----------------------
$excel.ActiveWindow.WindowState = -4140  # xlMinimized
